$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '93.772.15'
$ws.Range('E2').Value = '  -4.33%  '
$ws.Range('D3').Value = '3.397.99'
$ws.Range('E3').Value = '  -0.13%  '
$ws.Range('E4').Value = '  +0.08%  '
$ws.Range('D5').Value = '235.91'
$ws.Range('E5').Value = '  -7.00%  '
$ws.Range('D6').Value = '636.66'
$ws.Range('E6').Value = '  -3.81%  '
$ws.Range('E7').Value = '  -3.69%  '
$ws.Range('D8').Value = '0.398'
$ws.Range('E8').Value = '  -7.07%  '
$ws.Range('E9').Value = '  +0.15%  '
$ws.Range('D10').Value = '0.958'
$ws.Range('E10').Value = '  -7.76%  '
$ws.Range('D11').Value = '3.396.32'
$ws.Range('E12').Value = '  -5.31%  '
$ws.Range('D13').Value = '41.10'
$ws.Range('E13').Value = '  -7.94%  '
$ws.Range('D14').Value = '6.14'
$ws.Range('E14').Value = '  +0.35%  '
$ws.Range('D15').Value = '93.630.49'
$ws.Range('E15').Value = '  -4.22%  '
$ws.Range('D16').Value = '4.032.53'
$ws.Range('E16').Value = '  -0.51%  '
$ws.Range('D17').Value = '0.0000247'
$ws.Range('E17').Value = '  -4.33%  '
$ws.Range('D18').Value = '8.21'
$ws.Range('E18').Value = '  -10.69%  '
$ws.Range('D19').Value = '3.395.75'
$ws.Range('E19').Value = '  +0.36%  '
$ws.Range('D20').Value = '17.24'
$ws.Range('E20').Value = '  -5.39%  '
$ws.Range('D21').Value = '11.42'
$ws.Range('E21').Value = '  -0.89%  '
$ws.Range('D22').Value = '495.17'
$ws.Range('E22').Value = '  -3.24%  '
$ws.Range('D23').Value = '0.463'
$ws.Range('E23').Value = '  -11.47%  '
$ws.Range('D24').Value = '3.22'
$ws.Range('E24').Value = '  -5.95%  '
$ws.Range('D25').Value = '0.0000189'
$ws.Range('E25').Value = '  -6.27%  '
$ws.Range('D26').Value = '6.42'
$ws.Range('E26').Value = '  -6.93%  '
$ws.Range('D27').Value = '90.67'
$ws.Range('E27').Value = '  -6.79%  '
$ws.Range('D28').Value = '3.581.22'
$ws.Range('E28').Value = '  -0.10%  '
$ws.Range('D29').Value = '11.75'
$ws.Range('E29').Value = '  -5.82%  '
$ws.Range('D30').Value = '11.37'
$ws.Range('E30').Value = '  -5.08%  '
$ws.Range('E31').Value = '  -0.14%  '
$ws.Range('D32').Value = '2.70'
$ws.Range('E32').Value = '  -0.89%  '
$ws.Range('E33').Value = '  -6.12%  '
$ws.Range('E34').Value = '  +0.35%  '
$ws.Range('E35').Value = '  -8.47%  '
$ws.Range('D36').Value = '29.15'
$ws.Range('E36').Value = '  +0.10%  '
$ws.Range('D37').Value = '0.540'
$ws.Range('E37').Value = '  -4.31%  '
$ws.Range('D38').Value = '532.66'
$ws.Range('E38').Value = '  +1.43%  '
$ws.Range('D39').Value = '7.58'
$ws.Range('E39').Value = '  -4.61%  '
$ws.Range('D40').Value = '1.42'
$ws.Range('E40').Value = '  -4.73%  '
$ws.Range('E42').Value = '  -2.97%  '
$ws.Range('D43').Value = '0.891'
$ws.Range('E43').Value = '  +3.32%  '
$ws.Range('D44').Value = '24.01'
$ws.Range('E44').Value = '  -1.64%  '
$ws.Range('D45').Value = '3.69'
$ws.Range('E45').Value = '  -0.03%  '
$ws.Range('E46').Value = '  -2.90%  '
$ws.Range('D47').Value = '5.57'
$ws.Range('E47').Value = '  -0.98%  '
$ws.Range('D48').Value = '2.17'
$ws.Range('E48').Value = '  -3.20%  '
$ws.Range('B49').Value = 'OKB'
$ws.Range('C49').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D49').Value = '54.20'
$ws.Range('E49').Value = '  -3.55%  '
$ws.Range('B50').Value = 'dogwifhat'
$ws.Range('C50').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D50').Value = '3.24'
$ws.Range('E50').Value = '  +0.07%  '
$ws.Range('D51').Value = '0.0399'
$ws.Range('E51').Value = '  -6.26%  '
